$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.432.57'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.63%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.968.57'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -5.38%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.79'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.09%  '
$ws.Range("E6").Value = '  -4.59%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.57'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -9.59%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").Value = '  -5.51%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '55.78'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.59%  '
$ws.Range("E11").Value = '  +3.85%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.103'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.00'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -6.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.832'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -10.61%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.256.60'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.42%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.53'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -8.68%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.34'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.39%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.970.22'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -5.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '36.357.41'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.76%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.11'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.54%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0888'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.37%  '
$ws.Range("B22").Value = 'BitcoinCash'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '231.50'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.59%  '
$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.12'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -7.07%  '
$ws.Range("E24").Value = '  -0.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.54'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.39%  '
$ws.Range("E26").Value = '  -4.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.62'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '164.10'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.35%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.86'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.42%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.124'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.71%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.118'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.87%  '
$ws.Range("E32").Value = '  -3.21%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.75'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -8.17%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0639'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.82%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.33'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -7.39%  '
$ws.Range("E36").Value = '  +0.04%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.81'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.55%  '
$ws.Range("E38").Value = '  -7.14%  '
$ws.Range("E39").Value = '  -15.34%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.91'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.45%  '
$ws.Range("E41").Value = '  -6.04%  '
$ws.Range("E42").Value = '  -4.70%  '
$ws.Range("E43").Value = '  -8.33%  '
$ws.Range("E44").Value = '  -3.96%  '
$ws.Range("E45").Value = '  -9.77%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '15.80'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -8.38%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '88.80'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.58%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.348.96'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.80%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.27'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -9.12%  '
$ws.Range("E50").Value = '  -4.22%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '44.67'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.24%  '
